$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$sub3 = [char]8323

$dcell = $ws.Range("D2")
$dcell.NumberFormat = "@"
$dcell.Value = "36.922.12"
$dcell.Style = "Normal"
$ws.Range("E2").Value = "  -1.56%  "

$dcell = $ws.Range("D3")
$dcell.NumberFormat = "@"
$dcell.Value = "2.017.28"
$dcell.Style = "Normal"
$ws.Range("E3").Value = "  -3.00%  "

$dcell = $ws.Range("D4")
$dcell.NumberFormat = "@"
$dcell.Value = "0.997"
$dcell.Style = "Normal"
$ws.Range("E4").Value = "  -0.31%  "

$dcell = $ws.Range("D5")
$dcell.NumberFormat = "@"
$dcell.Value = "226.42"
$dcell.Style = "Normal"
$ws.Range("E5").Value = "  -2.71%  "

$dcell = $ws.Range("D6")
$dcell.NumberFormat = "@"
$dcell.Value = "0.607"
$dcell.Style = "Normal"
$ws.Range("E6").Value = "  -4.14%  "

$ws.Range("E7").Value = "  +0.04%  "

$dcell = $ws.Range("D8")
$dcell.NumberFormat = "@"
$dcell.Value = "54.87"
$dcell.Style = "Normal"
$ws.Range("E8").Value = "  -5.05%  "

$dcell = $ws.Range("D9")
$dcell.NumberFormat = "@"
$dcell.Value = "0.379"
$dcell.Style = "Normal"
$ws.Range("E9").Value = "  -2.68%  "

$ws.Range("E10").Value = "  +1.56%  "

$dcell = $ws.Range("D11")
$dcell.NumberFormat = "@"
$dcell.Value = "0.105"
$dcell.Style = "Normal"
$ws.Range("E11").Value = "  -3.55%  "

$dcell = $ws.Range("D12")
$dcell.NumberFormat = "@"
$dcell.Value = "2.314.63"
$dcell.Style = "Normal"
$ws.Range("E12").Value = "  -3.03%  "

$dcell = $ws.Range("D13")
$dcell.NumberFormat = "@"
$dcell.Value = "14.27"
$dcell.Style = "Normal"
$ws.Range("E13").Value = "  -5.20%  "

$dcell = $ws.Range("D14")
$dcell.NumberFormat = "@"
$dcell.Value = "20.52"
$dcell.Style = "Normal"
$ws.Range("E14").Value = "  -2.42%  "

$ws.Range("E15").Value = "  -3.47%  "

$ws.Range("E16").Value = "  -3.42%  "

$dcell = $ws.Range("D17")
$dcell.NumberFormat = "@"
$dcell.Value = "2.020.47"
$dcell.Style = "Normal"
$ws.Range("E17").Value = "  -2.67%  "

$dcell = $ws.Range("D18")
$dcell.NumberFormat = "@"
$dcell.Value = "36.820.10"
$dcell.Style = "Normal"
$ws.Range("E18").Value = "  -1.66%  "

$dcell = $ws.Range("D19")
$dcell.NumberFormat = "@"
$dcell.Value = "6.10"
$dcell.Style = "Normal"
$ws.Range("E19").Value = "  +1.04%  "

$dcell = $ws.Range("D20")
$dcell.NumberFormat = "@"
$dcell.Value = "68.81"
$dcell.Style = "Normal"
$ws.Range("E20").Value = "  -2.65%  "

$dcell = $ws.Range("D21")
$dcell.NumberFormat = "@"
$dcell.Value = "0.0${sub3}0826"
$dcell.Style = "Normal"
$ws.Range("E21").Value = "  -0.54%  "

$dcell = $ws.Range("D22")
$dcell.NumberFormat = "@"
$dcell.Value = "226.46"
$dcell.Style = "Normal"
$ws.Range("E22").Value = "  -0.97%  "

$ws.Range("E24").Value = "  +2.86%  "

$dcell = $ws.Range("D25")
$dcell.NumberFormat = "@"
$dcell.Value = "2.26"
$dcell.Style = "Normal"
$ws.Range("E25").Value = "  -4.76%  "

$dcell = $ws.Range("D26")
$dcell.NumberFormat = "@"
$dcell.Value = "167.23"
$dcell.Style = "Normal"
$ws.Range("E26").Value = "  -1.76%  "

$dcell = $ws.Range("D27")
$dcell.NumberFormat = "@"
$dcell.Value = "9.24"
$dcell.Style = "Normal"
$ws.Range("E27").Value = "  -5.04%  "

$dcell = $ws.Range("D28")
$dcell.NumberFormat = "@"
$dcell.Value = "0.126"
$dcell.Style = "Normal"
$ws.Range("E28").Value = "  -5.31%  "

$ws.Range("E29").Value = "  -3.95%  "

$dcell = $ws.Range("D30")
$dcell.NumberFormat = "@"
$dcell.Value = "1.33"
$dcell.Style = "Normal"
$ws.Range("E30").Value = "  -3.25%  "

$ws.Range("E31").Value = "  -4.56%  "

$dcell = $ws.Range("D32")
$dcell.NumberFormat = "@"
$dcell.Value = "4.48"
$dcell.Style = "Normal"
$ws.Range("E32").Value = "  -3.67%  "

$dcell = $ws.Range("D33")
$dcell.NumberFormat = "@"
$dcell.Value = "0.0613"
$dcell.Style = "Normal"
$ws.Range("E33").Value = "  -3.43%  "

$ws.Range("E34").Value = "  -4.50%  "

$ws.Range("E35").Value = "  -4.01%  "

$ws.Range("E36").Value = "  +0.75%  "

$ws.Range("E37").Value = "  +0.00%  "

$dcell = $ws.Range("D38")
$dcell.NumberFormat = "@"
$dcell.Value = "3.17"
$dcell.Style = "Normal"
$ws.Range("E38").Value = "  -4.56%  "

$dcell = $ws.Range("D39")
$dcell.NumberFormat = "@"
$dcell.Value = "5.37"
$dcell.Style = "Normal"
$ws.Range("E39").Value = "  +1.01%  "

$ws.Range("E40").Value = "  -5.25%  "

$dcell = $ws.Range("D41")
$dcell.NumberFormat = "@"
$dcell.Value = "1.488.20"
$dcell.Style = "Normal"
$ws.Range("E41").Value = "  +2.16%  "

$dcell = $ws.Range("D42")
$dcell.NumberFormat = "@"
$dcell.Value = "16.97"
$dcell.Style = "Normal"
$ws.Range("E42").Value = "  +0.80%  "

$ws.Range("E43").Value = "  -2.66%  "

$dcell = $ws.Range("D44")
$dcell.NumberFormat = "@"
$dcell.Value = "95.06"
$dcell.Style = "Normal"
$ws.Range("E44").Value = "  -5.14%  "

$dcell = $ws.Range("D45")
$dcell.NumberFormat = "@"
$dcell.Value = "2.78"
$dcell.Style = "Normal"
$ws.Range("E45").Value = "  -4.38%  "

$dcell = $ws.Range("D46")
$dcell.NumberFormat = "@"
$dcell.Value = "1.14"
$dcell.Style = "Normal"
$ws.Range("E46").Value = "  -4.85%  "

$dcell = $ws.Range("D47")
$dcell.NumberFormat = "@"
$dcell.Value = "7.29"
$dcell.Style = "Normal"
$ws.Range("E47").Value = "  +0.38%  "

$ws.Range("E48").Value = "  -4.35%  "

$ws.Range("E49").Value = "  -1.35%  "

$dcell = $ws.Range("D50")
$dcell.NumberFormat = "@"
$dcell.Value = "3.67"
$dcell.Style = "Normal"
$ws.Range("E50").Value = "  -7.53%  "

$dcell = $ws.Range("D51")
$dcell.NumberFormat = "@"
$dcell.Value = "2.203.83"
$dcell.Style = "Normal"
$ws.Range("E51").Value = "  -2.97%  "
